# "Abstract, introduction, fazit bearbeitet"
#
# Adds a third ("combined") weighted-mean evaluation based on rows 9 and 21
# (N27/O27/P27), a new "uc" label (P26) driving that block's header row, a
# small helper conversion in D29, and renormalises the "Verhaeltnis der
# Peakflaechen" ratios in column V (rows 15-17) against B18 instead of B24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Renormalise the peak-area ratio column (V15:V17) against B18 ---
$ws.Range("V15").Formula = "=B18/B18"
$ws.Range("V16").Formula = "=B21/B18"
$ws.Range("V17").Formula = "=B24/B18"

# --- New "uc" header cell above the combined block ---
$ws.Range("P26").Value = "uc"

# --- New combined row (weighted mean across the two peak pairs) ---
$ws.Range("N27").Formula = "=AVERAGE(L21,L9)"
$ws.Range("O27").Formula = "=STDEV(L9,L21)/SQRT(2)"
$ws.Range("P27").Formula = "=SQRT((M9)^2+(M21)^2)"

# --- Small standalone helper conversion ---
$ws.Range("D29").Formula = "=110/(2*PI())"

# --- Final selection / scroll position as left by the author ---
$ws.Range("M28").Select()
